$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 17:33"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6466498
$ws.Range("C4").Value = 6248
$ws.Range("D4").Value = 3726119
$ws.Range("E4").Value = 2547096
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 193283

# Row 5: India
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 4236961
$ws.Range("C5").Value = 34399
$ws.Range("D5").Value = 3278999
$ws.Range("E5").Value = 885929
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 346
$ws.Range("H5").Value = 72033

# Row 14: Chile
$ws.Range("A14").Value = "Chile"
$ws.Range("B14").Value = 424274
$ws.Range("C14").Value = 1764
$ws.Range("D14").Value = 395717
$ws.Range("E14").Value = 16905
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 60
$ws.Range("H14").Value = 11652

# Row 24: Alemania
$ws.Range("A24").Value = "Alemania"
$ws.Range("B24").Value = 252714
$ws.Range("C24").Value = 990
$ws.Range("D24").Value = 227000
$ws.Range("E24").Value = 16313
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 9401

# Row 31: Catar
$ws.Range("A31").Value = "Catar"
$ws.Range("B31").Value = 120348
$ws.Range("C31").Value = 253
$ws.Range("D31").Value = 117241
$ws.Range("E31").Value = 2902
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 205

# Row 53: Singapur
$ws.Range("A53").Value = "Singapur"
$ws.Range("B53").Value = 57044
$ws.Range("C53").Value = 22
$ws.Range("D53").Value = 56408
$ws.Range("E53").Value = 609
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 27

# Row 65: Moldavia
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 40055
$ws.Range("C65").Value = 258
$ws.Range("D65").Value = 27799
$ws.Range("E65").Value = 11182
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 11
$ws.Range("H65").Value = 1074

# Row 68: Kenia
$ws.Range("A68").Value = "Kenia"
$ws.Range("B68").Value = 35205
$ws.Range("C68").Value = 102
$ws.Range("D68").Value = 21310
$ws.Range("E68").Value = 13296
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 599

# Row 93: Albania
$ws.Range("A93").Value = "Albania"
$ws.Range("B93").Value = 10406
$ws.Range("C93").Value = 151
$ws.Range("D93").Value = 6186
$ws.Range("E93").Value = 3901
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = 319

# Row 97: Guayana Francesa
$ws.Range("A97").Value = "Guayana Francesa"
$ws.Range("B97").Value = 9355
$ws.Range("C97").Value = 33
$ws.Range("D97").Value = 8902
$ws.Range("E97").Value = 391
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 62

# Row 100: Namibia
$ws.Range("A100").Value = "Namibia"
$ws.Range("B100").Value = 8810
$ws.Range("C100").Value = 125
$ws.Range("D100").Value = 3806
$ws.Range("E100").Value = 4913
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = 91

# Row 124: Uganda
$ws.Range("A124").Value = "Uganda"
$ws.Range("B124").Value = 3776
$ws.Range("C124").Value = 109
$ws.Range("D124").Value = 1741
$ws.Range("E124").Value = 1991
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 3
$ws.Range("H124").Value = 44

# Row 141: Reunion
$ws.Range("A141").Value = "Reunion"
$ws.Range("B141").Value = 2277
$ws.Range("C141").Value = 55
$ws.Range("D141").Value = 1313
$ws.Range("E141").Value = 951
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 13

# Row 142: Trinidad yTobago
$ws.Range("A142").Value = "Trinidad yTobago"
$ws.Range("B142").Value = 2250
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 724
$ws.Range("E142").Value = 1492
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 34

# Row 143: Guinea-Bisau
$ws.Range("A143").Value = "Guinea-Bisau"
$ws.Range("B143").Value = 2245
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 1127
$ws.Range("E143").Value = 1080
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 38

# Row 195: Curazao
$ws.Range("A195").Value = "Curazao"
$ws.Range("B195").Value = 102
$ws.Range("C195").Value = 10
$ws.Range("D195").Value = 45
$ws.Range("E195").Value = 56
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 1

# Row 196: Antigua y Barbuda
$ws.Range("A196").Value = "Antigua y Barbuda"
$ws.Range("B196").Value = 95
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 91
$ws.Range("E196").Value = 1
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 3
